$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table contents (CON1 / CON2 / Coculture comparison rows expanded
# from 6 data rows to 10 data rows, with different fungal isolate codes).
$data = @(
    @("CON1",   "CON2",  "Coculture"),
    @("F3CON",  "F15CON","F3vF15"),
    @("F11CON", "F4CON", "F11vF4"),
    @("F11CON", "F2CON", "F11vF2"),
    @("F6CON",  "F2CON", "F6vF2"),
    @("F15CON", "F9CON", "F15vF9"),
    @("F14CON", "F13CON","F14vF13"),
    @("F1CON",  "F6CON", "F1vF6"),
    @("F9CON",  "F7CON", "F9vF7"),
    @("F5CON",  "F6CON", "F5vF6"),
    @("F14CON", "F8CON", "F14vF8")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Switch the page to landscape / A4 like the resaved workbook.
$ws.PageSetup.Orientation = 2
$ws.PageSetup.PaperSize = 9

# Move the active selection like in the author's final save.
$ws.Range("D9").Select()
